# "Update Master dan lainnya" - remove the "mrp" and "valcl" columns from the
# material master sheet (originally columns C and G); remaining columns
# shift left so the table becomes material | deskripsi | mtyp | matl | bun.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column G (valcl) first, then column C (mrp); order doesn't matter
# here since the two columns don't interact, but this mirrors how a user
# would clean up the sheet right-to-left.
$ws.Range("G:G").EntireColumn.Delete()
$ws.Range("C:C").EntireColumn.Delete()

# Move the active cell/selection to C1 (matches the saved selection state).
$ws.Range("C1").Select()
